# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for rows 2-29 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(2, 3, 5, 10, 8, 2, 3, 3, 6, 7, 4, 7, 7, 9, 10, 6, 9, 4, 8, 2, 0, 3, 9, 5, 5, 1, 2, 1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
